# specs_mw_one_scenario.xlsx edit
# - rename "SA_PV_cost" / "SA_PV_Cost" header text to "PV_cost_adjust" / "PV_Cost_adjust"
# - update PV cost level values on ScenarioParameters (G2/G3: 4300->1, 5500->1.25)
# - remove the TimeStep column from SpecsData / SpecsDataCalib / SpecsDataCalib1
# - make ScenarioParameters the active tab, update various sheet selections

$wb = $excel.ActiveWorkbook

# --- ScenarioInfo -----------------------------------------------------
$wsInfo = $wb.Worksheets.Item("ScenarioInfo")
$wsInfo.Range("F1").Value = "PV_cost_adjust"
$wsInfo.Range("F1").Select()

# --- ScenarioParameters -------------------------------------------------
$wsParams = $wb.Worksheets.Item("ScenarioParameters")
$wsParams.Range("G1").Value = "PV_Cost_adjust"
$wsParams.Range("G2").Value = 1
$wsParams.Range("G3").Value = 1.25

# --- SpecsData / SpecsDataCalib / SpecsDataCalib1: drop TimeStep column (E) ---
$wsData = $wb.Worksheets.Item("SpecsData")
$wsData.Columns("E:E").Delete()
$wsData.Range("B4").Select()

$wsCalib = $wb.Worksheets.Item("SpecsDataCalib")
$wsCalib.Columns("E:E").Delete()
$wsCalib.Range("E:E").Select()

$wsCalib1 = $wb.Worksheets.Item("SpecsDataCalib1")
$wsCalib1.Columns("E:E").Delete()
$wsCalib1.Range("G11").Select()

# --- Activate ScenarioParameters tab + set its selection last ---------
$wsParams.Activate()
$wsParams.Range("H6").Select()
